$wb = $excel.ActiveWorkbook

# --- Rename sheets (sheet tab names updated with new timestamps) ---
$wb.Sheets.Item(1).Name = "GNG_TO-16502911163442667"
$wb.Sheets.Item(2).Name = "NB_TO-16502911187393475"
$wb.Sheets.Item(3).Name = "RS_TO-16502911187413082"
$wb.Sheets.Item(4).Name = "TOL_TO-1650291118788176"
$wb.Sheets.Item(5).Name = "vSAT_TO-1650291118866368"

# --- Sheet 1: GNG ---
$ws1 = $wb.Sheets.Item(1)
$ws1.Range("B2").Value = "go_stims-1650291116307487.csv"
$ws1.Range("B3").Value = "GNG_stims-16502911163283188.csv"
$ws1.Range("B4").Value = "go_stims-1650291116329327.csv"
$ws1.Range("B5").Value = "GNG_stims-16502911163433118.csv"

# --- Sheet 2: NB ---
$ws2 = $wb.Sheets.Item(2)
$ws2.Range("B2").Value = "ZB-match_5-16502911164968307.csv"
$ws2.Range("B3").Value = "OB-16502911167638822.csv"
$ws2.Range("B4").Value = "OB-16502911170757525.csv"
$ws2.Range("B5").Value = "TB-16502911180517466.csv"
$ws2.Range("B6").Value = "ZB-match_0-16502911167077854.csv"
$ws2.Range("B7").Value = "TB-1650291118727186.csv"
$ws2.Range("B8").Value = "TB-16502911171515625.csv"
$ws2.Range("B9").Value = "ZB-match_8-16502911164671624.csv"
$ws2.Range("B10").Value = "OB-16502911169699943.csv"

# --- Sheet 3: RS (no cell content changes) ---

# --- Sheet 4: TOL ---
$ws4 = $wb.Sheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16502911187553182.csv"
$ws4.Range("B3").Value = "ZM_stims-16502911187423563.csv"
$ws4.Range("B4").Value = "MM_stims-16502911187718277.csv"
$ws4.Range("B5").Value = "ZM_stims-16502911187553182.csv"
$ws4.Range("B6").Value = "MM_stims-16502911187872167.csv"
$ws4.Range("B7").Value = "ZM_stims-16502911187728293.csv"

# --- Sheet 5: vSAT ---
$ws5 = $wb.Sheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-1650291118818341.csv"
$ws5.Range("B3").Value = "SAT_stims-16502911187920463.csv"
$ws5.Range("B4").Value = "vSAT_stims-1650291118850211.csv"
$ws5.Range("B5").Value = "vSAT_stims-16502911188346283.csv"
